$wb = $excel.ActiveWorkbook

# --- Sheet "tr0001": give B1/B2 the same (non-default) cell style already used
# by their row-mates (A1/C1, A2/C2). A format round-trip (change then revert)
# nudges the engine into assigning an explicit style record instead of the
# implicit default one.
$ws1 = $wb.Worksheets.Item("tr0001")
$ws1.Range("B1").Font.Name = "Calibri"
$ws1.Range("B1").Font.Name = "Arial"
$ws1.Range("B2").Font.Name = "Calibri"
$ws1.Range("B2").Font.Name = "Arial"

# --- Sheet "ts0001": split the single comma-joined "keywords" cell (D1) into
# one keyword per row (D1, D2, D3).
$ws2 = $wb.Worksheets.Item("ts0001")
$ws2.Range("D1").Value = "test"
$ws2.Range("D2").Value = "test data"
$ws2.Range("D3").Value = "unkown results"

# --- Switch the active sheet/tab from "tr0001" to "ts0001" and move that
# sheet's selection onto the newly-filled D2 cell.
$ws2.Activate()
$ws2.Range("D2").Select()
